# dice_cov_fitness.xlsx edit
# "finished test_fitness.m, nearly done with test_substances.m"
#
#  - rename the still-default "Sheet11" to "sitting_7days_mm"
#  - move on to "influ_exer_ny" (becomes the active/selected tab), leaving
#    a new selection there
#  - leave behind an updated selection on "last7_newactivities_ny" too
#  - nudge the app window around (best effort; harmless if the host
#    doesn't persist geometry)

$wb = $excel.ActiveWorkbook

# 1. Rename the placeholder sheet to its real name.
$sittingSheet = $wb.Worksheets.Item("Sheet11")
$sittingSheet.Name = "sitting_7days_mm"

# 2. Update the lingering selection on last7_newactivities_ny (A24 -> E55),
#    without leaving it as the active tab.
$last7Sheet = $wb.Worksheets.Item("last7_newactivities_ny")
$last7Sheet.Activate()
$last7Sheet.Range("E55").Select()

# 3. Move work to influ_exer_ny: it becomes the active/selected sheet, with
#    a fresh single-cell selection (S47) replacing the old A16:L20 block.
$influSheet = $wb.Worksheets.Item("influ_exer_ny")
$influSheet.Activate()
$influSheet.Range("S47").Select()

# 4. Reposition/resize the workbook window to match the author's session.
$win = $excel.ActiveWindow
$win.Left = 38270
$win.Top = 70
$win.Width = 38660
$win.Height = 21260
